$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bid-opening-date column (J) currently holds text like "25.04.03".
# Replace with real date serials and apply a date number format (Excel
# normalizes "mm-dd-yy" to the built-in numFmtId 14 short-date format).
$dates = @{
    2  = 45750
    3  = 45744
    4  = 45742
    5  = 45736
    6  = 45735
    7  = 45734
    8  = 45733
    9  = 45707
    10 = 45679
    11 = 45366
    12 = 45359
    13 = 45356
    14 = 45350
    15 = 45345
    16 = 45337
    17 = 45337
    18 = 45321
    19 = 45308
    20 = 45307
    21 = 45300
    22 = 45231
    23 = 44985
    24 = 44985
    25 = 44980
    26 = 44980
    27 = 44644
    28 = 44630
    29 = 44630
    30 = 44610
    31 = 44588
    32 = 44554
}

# Set the format on the first cell, then fan that single new style out to
# the rest of the column via a format-only paste so every row shares one
# cellXf (matches how Excel itself reuses styles for a column-wide edit).
$first = $ws.Cells.Item(2, 10)
$first.Value = $dates[2]
$first.NumberFormat = "mm-dd-yy"

$restRange = $ws.Range("J3:J32")
$first.Copy()
$restRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($row = 3; $row -le 32; $row++) {
    $ws.Cells.Item($row, 10).Value = $dates[$row]
}

# Match the author's resulting view state (scrolled down, B19 selected).
$ws.Range("B19").Select()
$excel.ActiveWindow.ScrollRow = 11
